$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.083.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +6.24%  "
$ws.Range("D3").Value = "'1.716.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.73%  "
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").Value = "'333.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.48%  "
$ws.Range("D6").Value = "'0.9993"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").Value = "'0.3693"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.63%  "
$ws.Range("D8").Value = "'49.31"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.38%  "
$ws.Range("D9").Value = "'0.3343"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'1.186"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.75%  "
$ws.Range("D11").Value = "'0.07464"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.98%  "
$ws.Range("D12").Value = "'0.9996"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("D13").Value = "'6.296"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.42%  "
$ws.Range("D14").Value = "'20.06"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.92%  "
$ws.Range("D15").Value = "'6.936"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.87%  "
$ws.Range("D16").Value = "'1.714.21"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.20%  "
$ws.Range("D17").Value = "'0.00001077"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.24%  "
$ws.Range("D18").Value = "'0.06637"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("E19").Value = "  +3.92%  "
$ws.Range("D20").Value = "'0.9996"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").Value = "'16.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.35%  "
$ws.Range("D22").Value = "'6.086"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.77%  "
$ws.Range("D23").Value = "'13.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.90%  "
$ws.Range("D24").Value = "'26.032.62"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.89%  "
$ws.Range("E25").Value = "  +0.89%  "
$ws.Range("D26").Value = "'2.456"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.16%  "
$ws.Range("D27").Value = "'150.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.02%  "
$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D28").Value = "'1.342"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +11.40%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'19.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.25%  "
$ws.Range("D30").Value = "'1.902.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.19%  "
$ws.Range("D31").Value = "'129.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.84%  "
$ws.Range("D32").Value = "'4.109"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.12%  "
$ws.Range("D33").Value = "'5.930"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.72%  "
$ws.Range("D34").Value = "'1.721"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.86%  "
$ws.Range("D35").Value = "'0.08515"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.60%  "
$ws.Range("D36").Value = "'12.88"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.00%  "
$ws.Range("D37").Value = "'5.353"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.81%  "
$ws.Range("D38").Value = "'0.06237"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.58%  "
$ws.Range("D39").Value = "'0.02294"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.91%  "
$ws.Range("D40").Value = "'0.2139"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.23%  "
$ws.Range("D41").Value = "'8.559"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.96%  "
$ws.Range("D42").Value = "'1.225"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.97%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.6163"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.10%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'14.41"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +12.78%  "
$ws.Range("D45").Value = "'1.000"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").Value = "'3.831"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.34%  "
$ws.Range("D47").Value = "'0.5893"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.79%  "
$ws.Range("D48").Value = "'128.44"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.02%  "
$ws.Range("D49").Value = "'2.020"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.56%  "
$ws.Range("D50").Value = "'0.07264"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.47%  "
$ws.Range("D51").Value = "'77.03"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.48%  "
